$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.390.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.312.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.306.72'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.622'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.114'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '39.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.814.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.15'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.347.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.186.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.992'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("E21").Value = '  +5.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '296.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.172'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.113'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.74%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +16.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0491'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.996'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.292'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '134.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.22%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.131.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.620.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.34%  '
